# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (want-to-go count) figures in column F
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types)
# sheets. 本地生活 has no data rows so it is untouched.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 64
$wsExhibition.Range("F5").Value  = 231
$wsExhibition.Range("F8").Value  = 259
$wsExhibition.Range("F12").Value = 99
$wsExhibition.Range("F13").Value = 2211
$wsExhibition.Range("F14").Value = 56
$wsExhibition.Range("F15").Value = 32
$wsExhibition.Range("F16").Value = 509
$wsExhibition.Range("F17").Value = 509
$wsExhibition.Range("F22").Value = 1700
$wsExhibition.Range("F23").Value = 3853
$wsExhibition.Range("F25").Value = 61
$wsExhibition.Range("F27").Value = 1149
$wsExhibition.Range("F28").Value = 214
$wsExhibition.Range("F29").Value = 2046
$wsExhibition.Range("F32").Value = 88
$wsExhibition.Range("F33").Value = 279
$wsExhibition.Range("F35").Value = 457
$wsExhibition.Range("F36").Value = 668
$wsExhibition.Range("F38").Value = 398

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 24

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value  = 64
$wsAllTypes.Range("F5").Value  = 231
$wsAllTypes.Range("F8").Value  = 259
$wsAllTypes.Range("F12").Value = 99
$wsAllTypes.Range("F13").Value = 2211
$wsAllTypes.Range("F14").Value = 56
$wsAllTypes.Range("F15").Value = 24
$wsAllTypes.Range("F16").Value = 32
$wsAllTypes.Range("F17").Value = 509
$wsAllTypes.Range("F18").Value = 509
$wsAllTypes.Range("F23").Value = 1700
$wsAllTypes.Range("F24").Value = 3853
$wsAllTypes.Range("F26").Value = 61
$wsAllTypes.Range("F28").Value = 1149
$wsAllTypes.Range("F29").Value = 214
$wsAllTypes.Range("F30").Value = 2046
$wsAllTypes.Range("F33").Value = 88
$wsAllTypes.Range("F34").Value = 279
$wsAllTypes.Range("F36").Value = 457
$wsAllTypes.Range("F37").Value = 668
$wsAllTypes.Range("F39").Value = 398
